$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 658.06665
$ws.Range("I2").Value = 816.9091
$ws.Range("J2").Value = 221.25
$ws.Range("K2").Value = 816.9091
$ws.Range("L2").Value = 221.25
$ws.Range("M2").Value = -703.9091
$ws.Range("N2").Value = -447.25
$ws.Range("H17").Value = 1787.8334
$ws.Range("J17").Value = 1836.7727
$ws.Range("L17").Value = 5510.3181
$ws.Range("N17").Value = -5846.3181
$ws.Range("H33").Value = 596.2308
$ws.Range("J33").Value = 1024
$ws.Range("L33").Value = 1024
$ws.Range("N33").Value = -1482
$ws.Range("H40").Value = 6466.5
$ws.Range("J40").Value = 8999.5
$ws.Range("L40").Value = 8999.5
$ws.Range("N40").Value = -9349.5
$ws.Range("H69").Value = 8341.666999999999
$ws.Range("I69").Value = 8341.666999999999
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 25025.001
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -24151.001
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 8341.666999999999
$ws.Range("I72").Value = 8341.666999999999
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 75075.003
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -70707.003
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 1768.8889
$ws.Range("J80").Value = 2120.0833
$ws.Range("L80").Value = 6360.249899999999
$ws.Range("N80").Value = -8356.249899999999
$ws.Range("H82").Value = 4449
$ws.Range("I82").Value = 4449
$ws.Range("K82").Value = 13347
$ws.Range("M82").Value = -12941
$ws.Range("H83").Value = 1768.8889
$ws.Range("J83").Value = 2120.0833
$ws.Range("L83").Value = 19080.7497
$ws.Range("N83").Value = -29064.7497
$ws.Range("H85").Value = 4449
$ws.Range("I85").Value = 4449
$ws.Range("K85").Value = 13347
$ws.Range("M85").Value = -11943
$ws.Range("H97").Value = 1183.5
$ws.Range("J97").Value = 1209.7142
$ws.Range("L97").Value = 3629.1426
$ws.Range("N97").Value = -4621.142599999999
$ws.Range("H111").Value = 3104.4285
$ws.Range("J111").Value = 8750
$ws.Range("L111").Value = 26250
$ws.Range("N111").Value = -32384
$ws.Range("H112").Value = 2689
$ws.Range("J112").Value = 2563.1836
$ws.Range("L112").Value = 7689.550799999999
$ws.Range("N112").Value = -9905.550799999999
$ws.Range("H132").Value = 206333.61
$ws.Range("I132").Value = 241218.62
$ws.Range("K132").Value = 723655.86
$ws.Range("M132").Value = -721125.86
$ws.Range("H135").Value = 5698.5557
$ws.Range("I135").Value = 2366.625
$ws.Range("J135").Value = 10545
$ws.Range("K135").Value = 21299.625
$ws.Range("L135").Value = 94905
$ws.Range("M135").Value = -18764.625
$ws.Range("N135").Value = -99975

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4596.615
$ws.Range("I2").Value = 1892.1666
$ws.Range("K2").Value = 1892.1666
$ws.Range("M2").Value = -1779.1666
$ws.Range("H32").Value = 1853998.6
$ws.Range("I32").Value = 1738.0172
$ws.Range("K32").Value = 1738.0172
$ws.Range("M32").Value = -1451.0172
$ws.Range("H44").Value = 12345
$ws.Range("J44").Value = 12345
$ws.Range("L44").Value = 12345
$ws.Range("N44").Value = -13321
$ws.Range("H45").Value = 1824.3
$ws.Range("J45").Value = 2539.6
$ws.Range("L45").Value = 2539.6
$ws.Range("N45").Value = -3293.6
$ws.Range("H61").Value = 5509.2856
$ws.Range("I61").Value = 10284.857
$ws.Range("K61").Value = 10284.857
$ws.Range("M61").Value = -10072.857
$ws.Range("H80").Value = 12345
$ws.Range("J80").Value = 12345
$ws.Range("L80").Value = 12345
$ws.Range("N80").Value = -14341
$ws.Range("H83").Value = 12345
$ws.Range("J83").Value = 12345
$ws.Range("L83").Value = 37035
$ws.Range("N83").Value = -47019
$ws.Range("H102").Value = 2327.8696
$ws.Range("I102").Value = 2216.2856
$ws.Range("K102").Value = 2216.2856
$ws.Range("M102").Value = -594.2856000000002
$ws.Range("H110").Value = 5715.154
$ws.Range("I110").Value = 2722.25
$ws.Range("K110").Value = 2722.25
$ws.Range("M110").Value = -677.25
$ws.Range("H116").Value = 4596.615
$ws.Range("I116").Value = 1892.1666
$ws.Range("K116").Value = 1892.1666
$ws.Range("M116").Value = 401.8334
$ws.Range("H122").Value = 4510.7856
$ws.Range("I122").Value = 3929.25
$ws.Range("K122").Value = 11787.75
$ws.Range("M122").Value = -9337.75
$ws.Range("H132").Value = 683426.7
$ws.Range("I132").Value = 809038.7
$ws.Range("K132").Value = 2427116.1
$ws.Range("M132").Value = -2424586.1
$ws.Range("H136").Value = 5509.2856
$ws.Range("I136").Value = 10284.857
$ws.Range("K136").Value = 30854.571
$ws.Range("M136").Value = -28304.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4596.615
$ws.Range("I3").Value = 1892.1666
$ws.Range("K3").Value = 1892.1666
$ws.Range("M3").Value = -1778.1666
$ws.Range("H94").Value = 3743.5
$ws.Range("I94").Value = 2079.9
$ws.Range("K94").Value = 2079.9
$ws.Range("M94").Value = -1628.9
$ws.Range("H134").Value = 1390953.5
$ws.Range("I134").Value = 1819741.5
$ws.Range("K134").Value = 5459224.5
$ws.Range("M134").Value = -5456689.5
$ws.Range("H141").Value = 63026.9
$ws.Range("I141").Value = 69541.8
$ws.Range("J141").Value = 56512
$ws.Range("K141").Value = 69541.8
$ws.Range("L141").Value = 56512
$ws.Range("M141").Value = -64361.8
$ws.Range("N141").Value = -66872

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7057.647
$ws.Range("I62").Value = 7212.1816
$ws.Range("K62").Value = 7212.1816
$ws.Range("M62").Value = -6588.1816
$ws.Range("H65").Value = 7057.647
$ws.Range("I65").Value = 7212.1816
$ws.Range("K65").Value = 36060.908
$ws.Range("M65").Value = -32940.908
$ws.Range("H132").Value = 10200.833
$ws.Range("I132").Value = 5290.8945
$ws.Range("K132").Value = 15872.6835
$ws.Range("M132").Value = -13342.6835
$ws.Range("H134").Value = 43486696
$ws.Range("I134").Value = 58827536
$ws.Range("J134").Value = 20991.5
$ws.Range("K134").Value = 176482608
$ws.Range("L134").Value = 62974.5
$ws.Range("M134").Value = -176480073
$ws.Range("N134").Value = -68044.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1818349.9
$ws.Range("J12").Value = 221
$ws.Range("L12").Value = 663
$ws.Range("N12").Value = -1009
$ws.Range("H63").Value = 15999.571
$ws.Range("I63").Value = 5999
$ws.Range("K63").Value = 17997
$ws.Range("M63").Value = -17248
$ws.Range("H66").Value = 15999.571
$ws.Range("I66").Value = 5999
$ws.Range("K66").Value = 53991
$ws.Range("M66").Value = -50247
$ws.Range("H69").Value = 5281.5
$ws.Range("I69").Value = 3233.3333
$ws.Range("J69").Value = 7329.6665
$ws.Range("K69").Value = 9699.999899999999
$ws.Range("L69").Value = 21988.9995
$ws.Range("M69").Value = -8888.999899999999
$ws.Range("N69").Value = -23610.9995
$ws.Range("H72").Value = 5281.5
$ws.Range("I72").Value = 3233.3333
$ws.Range("J72").Value = 7329.6665
$ws.Range("K72").Value = 29099.9997
$ws.Range("L72").Value = 65966.9985
$ws.Range("M72").Value = -25043.9997
$ws.Range("N72").Value = -74078.9985
$ws.Range("H76").Value = 8015
$ws.Range("J76").Value = 8015
$ws.Range("L76").Value = 24045
$ws.Range("N76").Value = -24811
$ws.Range("H79").Value = 8015
$ws.Range("J79").Value = 8015
$ws.Range("L79").Value = 24045
$ws.Range("N79").Value = -26697
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3564
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8820
$ws.Range("N83").ClearContents()
$ws.Range("H87").Value = 23585.8
$ws.Range("I87").Value = 13965
$ws.Range("J87").Value = 29999.666
$ws.Range("K87").Value = 41895
$ws.Range("L87").Value = 89998.99800000001
$ws.Range("M87").Value = -40647
$ws.Range("N87").Value = -92494.99800000001
$ws.Range("H90").Value = 23585.8
$ws.Range("I90").Value = 13965
$ws.Range("J90").Value = 29999.666
$ws.Range("K90").Value = 125685
$ws.Range("L90").Value = 269996.994
$ws.Range("M90").Value = -119445
$ws.Range("N90").Value = -282476.994
$ws.Range("H129").Value = 2118.4
$ws.Range("I129").Value = 999.5
$ws.Range("J129").Value = 2398.125
$ws.Range("K129").Value = 2998.5
$ws.Range("L129").Value = 7194.375
$ws.Range("M129").Value = 2001.5
$ws.Range("N129").Value = -17194.375
$ws.Range("H136").Value = 45460640
$ws.Range("I136").Value = 20839712
$ws.Range("K136").Value = 62519136
$ws.Range("M136").Value = -62514036
$ws.Range("H140").Value = 37502650
$ws.Range("I140").Value = 50001800
$ws.Range("K140").Value = 150005400
$ws.Range("M140").Value = -150000220

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13323.583
$ws.Range("J70").Value = 9999
$ws.Range("L70").Value = 9999
$ws.Range("N70").Value = -10539
$ws.Range("H73").Value = 13323.583
$ws.Range("J73").Value = 9999
$ws.Range("L73").Value = 9999
$ws.Range("N73").Value = -11871
$ws.Range("H126").Value = 45470620
$ws.Range("I126").Value = 166672340
$ws.Range("J126").Value = 19974.625
$ws.Range("K126").Value = 500017020
$ws.Range("L126").Value = 59923.875
$ws.Range("M126").Value = -500014550
$ws.Range("N126").Value = -64863.875
$ws.Range("H132").Value = 5296.7393
$ws.Range("I132").Value = 5409.6113
$ws.Range("K132").Value = 16228.8339
$ws.Range("M132").Value = -13698.8339
$ws.Range("H136").Value = 21133.111
$ws.Range("J136").Value = 21133.111
$ws.Range("L136").Value = 63399.333
$ws.Range("N136").Value = -68499.333
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2121.697
$ws.Range("I82").Value = 1320
$ws.Range("J82").Value = 3725.0908
$ws.Range("K82").Value = 1320
$ws.Range("L82").Value = 3725.0908
$ws.Range("M82").Value = -959
$ws.Range("N82").Value = -4447.0908
$ws.Range("H85").Value = 2121.697
$ws.Range("I85").Value = 1320
$ws.Range("J85").Value = 3725.0908
$ws.Range("K85").Value = 1320
$ws.Range("L85").Value = 3725.0908
$ws.Range("M85").Value = -72
$ws.Range("N85").Value = -6221.0908
$ws.Range("H100").Value = 3913
$ws.Range("I100").Value = 6799.6
$ws.Range("J100").Value = 2469.7
$ws.Range("K100").Value = 6799.6
$ws.Range("L100").Value = 2469.7
$ws.Range("M100").Value = -6258.6
$ws.Range("N100").Value = -3551.7
$ws.Range("H132").Value = 10992279
$ws.Range("I132").Value = 20410390
$ws.Range("J132").Value = 4483.1665
$ws.Range("K132").Value = 61231170
$ws.Range("L132").Value = 13449.4995
$ws.Range("M132").Value = -61228640
$ws.Range("N132").Value = -18509.4995
$ws.Range("H136").Value = 52921670
$ws.Range("I136").Value = 111121864
$ws.Range("K136").Value = 333365592
$ws.Range("M136").Value = -333363042

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11720.417
$ws.Range("I62").Value = 11974.5
$ws.Range("J62").Value = 11593.375
$ws.Range("K62").Value = 11974.5
$ws.Range("L62").Value = 11593.375
$ws.Range("M62").Value = -11350.5
$ws.Range("N62").Value = -12841.375
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 11720.417
$ws.Range("I65").Value = 11974.5
$ws.Range("J65").Value = 11593.375
$ws.Range("K65").Value = 59872.5
$ws.Range("L65").Value = 57966.875
$ws.Range("M65").Value = -56752.5
$ws.Range("N65").Value = -64206.875
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 1924.8235
$ws.Range("I100").Value = 1531.0834
$ws.Range("K100").Value = 3062.1668
$ws.Range("M100").Value = -2521.1668
$ws.Range("H107").Value = 1647.7142
$ws.Range("I107").Value = 1339
$ws.Range("K107").Value = 4017
$ws.Range("M107").Value = -2097
$ws.Range("H110").Value = 42096
$ws.Range("J110").Value = 42096
$ws.Range("L110").Value = 42096
$ws.Range("N110").Value = -50276
$ws.Range("H132").Value = 3781631
$ws.Range("I132").Value = 4944.7896
$ws.Range("K132").Value = 14834.3688
$ws.Range("M132").Value = -12304.3688
$ws.Range("H136").Value = 23827108
$ws.Range("I136").Value = 38479172
$ws.Range("K136").Value = 115437516
$ws.Range("M136").Value = -115434966
